$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current values for columns D,H,I,J,K,L,M,N,O,P,Q (rows 2-24)
# before applying the permutation, since source and destination rows overlap.
$cols = @(4,8,9,10,11,12,13,14,15,16,17)
$snapshot = @{}
for ($r = 2; $r -le 24; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Mapping: destination row -> source row (original data to copy from)
$mapping = @{
    2 = 14
    3 = 24
    4 = 8
    5 = 11
    6 = 19
    7 = 17
    8 = 23
    9 = 10
    10 = 13
    11 = 7
    12 = 16
    13 = 18
    14 = 20
    15 = 21
    16 = 22
    17 = 3
    18 = 4
    19 = 12
    20 = 9
    21 = 5
    22 = 6
    23 = 2
    24 = 15
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
